# Generate Report for Handoff
# - zh-cn / de-de sheets: rows for 656fd5f4, 88560a67, 98c50234, 9c488297
#   change Priority from "low" to "ht" (handoff-triggered re-tag).
# - zh-cn sheet: Latest Handoff Datetime for those same rows moves from
#   2016-08-13 16:41:24 to 2016-08-13 16:41:49 (new handoff pass).
# - Overview / de-de sheets: Latest HO Xliff Generate Date / Latest
#   Handoff Datetime for the 656fd5f4 row moves from 16:41:33 to 16:41:57.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Priority column (E) on zh-cn and de-de: "low" -> "ht" for rows 4-7.
foreach ($r in 4..7) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
}

# Latest Handoff Datetime column (H) on zh-cn for rows 4-7: refreshed timestamp.
foreach ($r in 4..7) {
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-13 16:41:49"
}

# Latest HO Xliff Generate Date (Overview col G, rows 4-7) and the matching
# Latest Handoff Datetime on de-de (col H, rows 4-7) share the same string;
# update both so the underlying text moves from 16:41:33 to 16:41:57.
foreach ($r in 4..7) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-13 16:41:57"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-13 16:41:57"
}
